$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.037347197532654
$ws.Range("B1").Value = 0.9745578169822693
$ws.Range("C1").Value = 4.213390350341797
$ws.Range("D1").Value = 2.505427837371826
$ws.Range("E1").Value = 1.228139638900757
